$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.13 = 7632.62 pesos`n✅ 7632.62 pesos = 2.12 = 950.02 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update tasa values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 470
$ws2.Range("O10").Value = 3587.33
$ws2.Range("N12").Value = 3605
$ws2.Range("O12").Value = 448.71
